{"js": "// The underlying change is a pure re-serialization of <w:rPr> child-element\n// order (the bold/italic toggles move from immediately after <w:rStyle/>\n// to immediately after the font metrics, right before <w:color/>) \u2014 no\n// formatting value actually changes anywhere in the document. Re-applying\n// the existing Bold/Italic values through the object model is what makes\n// the host re-emit run/style properties in its normal element order, which\n// reproduces that reordering.\n\nconst body = context.document.body;\n\n// 1) Every run of text in the document body.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (const paragraph of paragraphs.items) {\n  paragraph.load(\"text\");\n  await context.sync();\n\n  const text = paragraph.text;\n  if (!text) {\n    continue;\n  }\n\n  // Search for the paragraph's own text to obtain a range scoped to just\n  // that run content (as opposed to paragraph.getRange(), which also spans\n  // the paragraph mark and would stamp an extra pPr/rPr on the paragraph).\n  const results = body.search(text, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const found of results.items) {\n    found.font.load(\"bold,italic\");\n    await context.sync();\n\n    found.font.bold = found.font.bold;\n    found.font.italic = found.font.italic;\n    await context.sync();\n  }\n}\n\n// 2) Every character style definition (e.g. \"Style\") used in the document.\nconst styles = context.document.getStyles();\nstyles.load(\"items\");\nawait context.sync();\n\nfor (const style of styles.items) {\n  style.load(\"type\");\n  await context.sync();\n\n  if (style.type !== \"Character\") {\n    continue;\n  }\n\n  style.font.load(\"bold,italic\");\n  await context.sync();\n\n  style.font.bold = style.font.bold;\n  style.font.italic = style.font.italic;\n  await context.sync();\n}\n", "ps1": "# The underlying change is a pure re-serialization of <w:rPr> child-element\n# order (the bold/italic toggles move from immediately after <w:rStyle/>\n# to immediately after the font metrics, right before <w:color/>) -- no\n# formatting value actually changes anywhere in the document. Re-applying\n# the existing Bold/Italic values through the object model is what makes\n# Word re-emit run/style properties in its normal element order, which\n# reproduces that reordering.\n\n$d = $word.ActiveDocument\n\n# 1) Every paragraph's run content in the document body. Use Start..(End-1)\n# rather than the paragraph's own Range so the paragraph mark itself is\n# excluded - touching the mark too would also stamp a pPr/rPr on the\n# paragraph, which the target change does not do.\nforeach ($p in $d.Paragraphs) {\n    $pStart = $p.Range.Start\n    $pEnd = $p.Range.End - 1\n    if ($pStart -lt $pEnd) {\n        $r = $d.Range($pStart, $pEnd)\n        $r.Font.Bold = $r.Font.Bold\n        $r.Font.Italic = $r.Font.Italic\n    }\n}\n\n# 2) Every character style definition (e.g. \"Style\") used in the document.\n$wdStyleTypeCharacter = 2\nforeach ($s in $d.Styles) {\n    if ($s.Type -eq $wdStyleTypeCharacter) {\n        $s.Font.Bold = $s.Font.Bold\n        $s.Font.Italic = $s.Font.Italic\n    }\n}\n"}
